# Generate Report for Handback
#
# The handback CI run processed a second file this time
# (eda8aa28-9073-4fd6-8907-f1bb0d8dfa47.md), in addition to refreshing the
# hashes/timestamps for the file that was already being tracked
# (18d60560-... is replaced by a newer handoff, c6644389-...).
#
# NOTE: values that look like booleans ("True"/"False") or empty strings are
# stored as literal text in this report (not native Boolean cells), so they
# are written with a leading apostrophe to force text interpretation - this
# mirrors how the existing cells in the workbook are already typed.

$wb = $excel.ActiveWorkbook

$oldGuid  = "18d60560-cda4-4d1e-bebc-f28c2cfcc9f6"
$newGuid1 = "c6644389-fa82-4909-a5dd-6302e36807e0"
$newGuid2 = "eda8aa28-9073-4fd6-8907-f1bb0d8dfa47"

$hash1 = "83fc92869e963fff0a75af38494750f9d7f26830"
$hash2 = "b8dabd8198b1f00b76614a727eff0d1809061028"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3")) | Out-Null

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "2016-08-28 17:01:35"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-28 17:01:35"
$ws.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/244099fffdd423e1fa575d94b2824041f77af1a0/e2e/$newGuid1.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid1.md"

$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/244099fffdd423e1fa575d94b2824041f77af1a0/e2e/$newGuid2.md", "", "", "e2e\$newGuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3")) | Out-Null

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.$hash1.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-28 17:01:30"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.$hash1.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-28 17:01:46"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = "$newGuid2.$hash2.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-28 17:01:30"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid2.$hash2.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-28 17:01:46"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/244099fffdd423e1fa575d94b2824041f77af1a0/e2e/$newGuid1.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "$newGuid1.md"
$ws.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/941ac0309e764de63b6526c3edc1dc269ae2acfc/e2e/$newGuid1.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "$newGuid1.md"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/244099fffdd423e1fa575d94b2824041f77af1a0/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/941ac0309e764de63b6526c3edc1dc269ae2acfc/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3")) | Out-Null

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.$hash1.de-de.xlf"
$ws.Range("H2").Value = "2016-08-28 17:01:35"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.$hash1.de-de.xlf"
$ws.Range("K2").Value = "2016-08-28 17:01:53"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = "$newGuid2.$hash2.de-de.xlf"
$ws.Range("H3").Value = "2016-08-28 17:01:35"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid2.$hash2.de-de.xlf"
$ws.Range("K3").Value = "2016-08-28 17:01:53"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/244099fffdd423e1fa575d94b2824041f77af1a0/e2e/$newGuid1.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "$newGuid1.md"
$ws.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/187a83c26c57703752b9b556b7ac0af96f348c7d/e2e/$newGuid1.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "$newGuid1.md"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/244099fffdd423e1fa575d94b2824041f77af1a0/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/187a83c26c57703752b9b556b7ac0af96f348c7d/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
